$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indicators")

# Update the dictionary (feature list) string in B2
$ws.Range("B2").Value = "LangScr, S1_DRP, S1_GRD_1TO19, S1_GRD_4TO49, S1_GRD_5TO59, S2_GRD_1TO19, S2_GRD_3TO39, S2_GRD_4TO49, S2_GRD_5TO59, S2_GRD_6TO7, S2_VS_S1, SchoolRegion_1, SchoolRegion_2, SchoolRegion_3, SchoolType_2, MotherEd_7, Campus_1, PostulationType_1"

# Update the recalculated metric values
$ws.Range("B4").Value = 0.77777777777777779
$ws.Range("B6").Value = 0.8125
$ws.Range("B7").Value = 0.8441558441558441

# Widen column B to fit the new (longer) dictionary string
$ws.Columns.Item(2).ColumnWidth = 247.333333
